# Update countries & provincias Spain
# Refresh the COVID dashboard data on sheet "Pais" to the later snapshot
# (27 Sep 2020, 11:06). Country rows keep their position/rank cell (A)
# pointing at the same "slot", but three countries (Eslovaquia, Lituania,
# Timor Oriental) jumped ahead of their neighbours in the ranking, and a
# number of rows received new totals for the later snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain data refreshes (ranking position unchanged) ---------------

# Filipinas (row 24)
$ws.Range("B24").Value = 304226
$ws.Range("C24").Value = 2995
$ws.Range("D24").Value = 252510
$ws.Range("E24").Value = 46372
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 5344

# Indonesia (row 26)
$ws.Range("B26").Value = 275213
$ws.Range("C26").Value = 3874
$ws.Range("D26").Value = 203014
$ws.Range("E26").Value = 61813
$ws.Range("G26").Value = 78
$ws.Range("H26").Value = 10386

# Israel (row 27)
$ws.Range("B27").Value = 229374
$ws.Range("C27").Value = 2274
$ws.Range("D27").Value = 159136
$ws.Range("E27").Value = 68788
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 1450

# Polonia (row 46)
$ws.Range("B46").Value = 87330
$ws.Range("C46").Value = 1350
$ws.Range("D46").Value = 67904
$ws.Range("E46").Value = 16994
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 2432

# Singapur (row 59)
$ws.Range("B59").Value = 57700
$ws.Range("C59").Value = 15
$ws.Range("E59").Value = 314

# Afganistan (row 69)
$ws.Range("B69").Value = 39227
$ws.Range("C69").Value = 35
$ws.Range("D69").Value = 32642
$ws.Range("E69").Value = 5132

# El Salvador (row 77)
$ws.Range("B77").Value = 28630
$ws.Range("C77").Value = 215
$ws.Range("D77").Value = 22879
$ws.Range("E77").Value = 4925

# Croacia (row 89)
$ws.Range("B89").Value = 16197
$ws.Range("C89").Value = 190
$ws.Range("D89").Value = 14609
$ws.Range("E89").Value = 1316
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 272

# Hong Kong (row 124)
$ws.Range("B124").Value = 5066
$ws.Range("C124").Value = 6
$ws.Range("D124").Value = 4786
$ws.Range("E124").Value = 175

# Letonia (row 161)
$ws.Range("B161").Value = 1676
$ws.Range("C161").Value = 22
$ws.Range("E161").Value = 336

# --- Eslovaquia overtakes Gabon and Haiti (rows 106-108) --------------

$ws.Range("A106").Value = "Eslovaquia"
$ws.Range("B106").Value = 9078
$ws.Range("C106").Value = 478
$ws.Range("D106").Value = 4178
$ws.Range("E106").Value = 4856
$ws.Range("H106").Value = 44

$ws.Range("A107").Value = "Gabon"
$ws.Range("B107").Value = 8728
$ws.Range("D107").Value = 7934
$ws.Range("E107").Value = 740
$ws.Range("H107").Value = 54

$ws.Range("A108").Value = "Haiti"
$ws.Range("B108").Value = 8723
$ws.Range("D108").Value = 6551
$ws.Range("E108").Value = 1945
$ws.Range("H108").Value = 227

# --- Lituania overtakes Trinidad y Tobago (rows 132-133) --------------

$ws.Range("A132").Value = "Lituania"
$ws.Range("B132").Value = 4385
$ws.Range("C132").Value = 90
$ws.Range("D132").Value = 2327
$ws.Range("E132").Value = 1967
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 91

$ws.Range("A133").Value = "Trinidad yTobago"
$ws.Range("B133").Value = 4312
$ws.Range("D133").Value = 2185
$ws.Range("E133").Value = 2057
$ws.Range("H133").Value = 70

# --- Timor Oriental overtakes Santa Lucia (rows 206-207) --------------
# (their totals happen to be identical this snapshot, so only the
# country names need to swap)

$ws.Range("A206").Value = "Timor Oriental"
$ws.Range("A207").Value = "Santa Lucia"

# --- Footer timestamp --------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 11:06"
